# "Generate Report for Handoff"
#
# The localization-status report moved from "In Translation" to
# "Ready for handoff": the Overview sheet and each per-locale sheet
# (zh-cn / de-de) get their Status text + the "Latest HO Xliff Generate
# Date" / "Latest Handoff Datetime" timestamp refreshed, and the Status
# column (now holding longer text) is widened to fit.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# --- Overview sheet: zh-cn / de-de status columns (E, F) + HO xliff date (G) ---
$wsOverview.Range("E2").Value = "Ready for handoff"
$wsOverview.Range("F2").Value = "Ready for handoff"
$wsOverview.Range("G2").Value = "2016-09-01 07:10:18"

# --- zh-cn sheet: Status (C) + Latest Handoff Datetime (H) ---
$wsZhCn.Range("C2").Value = "Ready for handoff"
$wsZhCn.Range("H2").Value = "2016-09-01 07:10:05"

# --- de-de sheet: Status (C) + Latest Handoff Datetime (H) ---
$wsDeDe.Range("C2").Value = "Ready for handoff"
$wsDeDe.Range("H2").Value = "2016-09-01 07:10:18"

# --- Widen the Status columns now that the text is longer ---
# (ColumnWidth is in characters; 16.25 is the closest reachable value to
# the authored target stored width of ~17.216 characters.)
$wsOverview.Columns.Item(5).ColumnWidth = 16.25
$wsOverview.Columns.Item(6).ColumnWidth = 16.25
$wsZhCn.Columns.Item(3).ColumnWidth = 16.25
$wsDeDe.Columns.Item(3).ColumnWidth = 16.25
